$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.228956228956229
$ws.Range("C2").Value = 0.5117845117845118
$ws.Range("J2").Value = 0.0101010101010101
$ws.Range("P2").Value = 0.1717171717171717
$ws.Range("S2").Value = 0.07744107744107744
$ws.Range("C3").Value = 0.04968944099378882
$ws.Range("J3").Value = 0.006211180124223602
$ws.Range("P3").Value = 0.7888198757763976
$ws.Range("S3").Value = 0.15527950310559
$ws.Range("J4").Value = 0.07272727272727272
$ws.Range("O4").Value = 0.01818181818181818
$ws.Range("P4").Value = 0.6727272727272727
$ws.Range("S4").Value = 0.2363636363636364
$ws.Range("B6").Value = 0.07281553398058252
$ws.Range("D6").Value = 0.009708737864077669
$ws.Range("F6").Value = 0.04368932038834952
$ws.Range("J6").Value = 0.2864077669902912
$ws.Range("O6").Value = 0.01941747572815534
$ws.Range("Q6").Value = 0.1747572815533981
$ws.Range("R6").Value = 0.09223300970873786
$ws.Range("S6").Value = 0.3009708737864077
$ws.Range("B7").Value = 0.1277777777777778
$ws.Range("D7").Value = 0.01666666666666667
$ws.Range("F7").Value = 0.02777777777777778
$ws.Range("J7").Value = 0.1833333333333333
$ws.Range("O7").Value = 0.02777777777777778
$ws.Range("Q7").Value = 0.1777777777777778
$ws.Range("R7").Value = 0.09444444444444444
$ws.Range("S7").Value = 0.3444444444444444
$ws.Range("B8").Value = 0.07795698924731183
$ws.Range("D8").Value = 0.02419354838709677
$ws.Range("F8").Value = 0.0564516129032258
$ws.Range("J8").Value = 0.1317204301075269
$ws.Range("O8").Value = 0.02150537634408602
$ws.Range("Q8").Value = 0.2043010752688172
$ws.Range("R8").Value = 0.1478494623655914
$ws.Range("S8").Value = 0.3360215053763441
$ws.Range("B9").Value = 0.1052631578947368
$ws.Range("D9").Value = 0.02631578947368421
$ws.Range("E9").Value = 0.006578947368421052
$ws.Range("F9").Value = 0.07236842105263158
$ws.Range("J9").Value = 0.1447368421052632
$ws.Range("O9").Value = 0.03289473684210526
$ws.Range("Q9").Value = 0.1907894736842105
$ws.Range("R9").Value = 0.1184210526315789
$ws.Range("S9").Value = 0.3026315789473684
$ws.Range("B10").Value = 0.1149881046788263
$ws.Range("D10").Value = 0.02934179222839017
$ws.Range("E10").Value = 0.0007930214115781126
$ws.Range("F10").Value = 0.07375099127676447
$ws.Range("J10").Value = 0.1371927042030135
$ws.Range("O10").Value = 0.01189532117367169
$ws.Range("Q10").Value = 0.2363203806502776
$ws.Range("R10").Value = 0.09199048374306107
$ws.Range("S10").Value = 0.3037272006344171
$ws.Range("G11").Value = 0.1118881118881119
$ws.Range("J11").Value = 0.1293706293706294
$ws.Range("K11").Value = 0.1853146853146853
$ws.Range("L11").Value = 0.5664335664335665
$ws.Range("S11").Value = 0.006993006993006993
$ws.Range("G12").Value = 0.7245508982035929
$ws.Range("J12").Value = 0.1976047904191617
$ws.Range("K12").Value = 0.01197604790419162
$ws.Range("L12").Value = 0.03592814371257485
$ws.Range("S12").Value = 0.02994011976047904
$ws.Range("G13").Value = 0.7317073170731707
$ws.Range("J13").Value = 0.1707317073170732
$ws.Range("S13").Value = 0.0975609756097561
$ws.Range("F15").Value = 0.015
$ws.Range("H15").Value = 0.12
$ws.Range("I15").Value = 0.08
$ws.Range("J15").Value = 0.38
$ws.Range("K15").Value = 0.065
$ws.Range("M15").Value = 0.01
$ws.Range("O15").Value = 0.06
$ws.Range("S15").Value = 0.27
$ws.Range("F16").Value = 0.009569377990430622
$ws.Range("H16").Value = 0.2105263157894737
$ws.Range("I16").Value = 0.05741626794258373
$ws.Range("J16").Value = 0.4114832535885167
$ws.Range("K16").Value = 0.1100478468899522
$ws.Range("M16").Value = 0.02392344497607655
$ws.Range("O16").Value = 0.01435406698564593
$ws.Range("S16").Value = 0.1626794258373206
$ws.Range("F17").Value = 0.01066098081023454
$ws.Range("H17").Value = 0.1769722814498934
$ws.Range("I17").Value = 0.09168443496801706
$ws.Range("J17").Value = 0.4136460554371002
$ws.Range("K17").Value = 0.1044776119402985
$ws.Range("M17").Value = 0.02345415778251599
$ws.Range("O17").Value = 0.07036247334754797
$ws.Range("S17").Value = 0.1087420042643923
$ws.Range("F18").Value = 0.03125
$ws.Range("H18").Value = 0.1785714285714286
$ws.Range("I18").Value = 0.05803571428571429
$ws.Range("J18").Value = 0.4553571428571428
$ws.Range("K18").Value = 0.1339285714285714
$ws.Range("M18").Value = 0.01339285714285714
$ws.Range("O18").Value = 0.03571428571428571
$ws.Range("S18").Value = 0.09375
$ws.Range("F19").Value = 0.01984126984126984
$ws.Range("H19").Value = 0.185515873015873
$ws.Range("I19").Value = 0.06746031746031746
$ws.Range("J19").Value = 0.3839285714285715
$ws.Range("K19").Value = 0.1150793650793651
$ws.Range("M19").Value = 0.02083333333333333
$ws.Range("O19").Value = 0.0873015873015873
$ws.Range("S19").Value = 0.1200396825396825

Write-Output "Updated 108 cells in Sheet1"
